$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Hoja1"

# Update header texts for FechaCreacion / FechaModificacion
$ws.Range("J1").Value = "FechaCreacion"
$ws.Range("K1").Value = "FechaModificacion"

# Apply the date-number-format style used by J2:K13 to the header cells too
$ws.Range("J2").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)

# Remove the autofilter (defined name _xlnm._FilterDatabase + worksheet autoFilter)
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# Remove frozen panes
$ws.Application.ActiveWindow.FreezePanes = $false

# Update selection to match new view (J1:K1048576 selected, active cell J1)
$ws.Range("J1:K1048576").Select()

$wb.Application.ActiveWindow.Zoom = 100
